$d = $word.ActiveDocument

# Update the date heading paragraph (2024-10-23 Wednesday -> 2024-10-24 Thursday)
$d.Paragraphs(1).Range.Text = "2024-10-24 Thursday"

# Update each math-problem cell in the table (20 rows x 5 columns, row-major order)
$table = $d.Tables(1)
$values = @(
    "71-38=",
    "20+70=",
    "97-15=",
    "85+5=",
    "87-23=",
    "30+5=",
    "59-48=",
    "70-40=",
    "0+70=",
    "57+1=",
    "78+20=",
    "81-13=",
    "28+17=",
    "91+4=",
    "14+26=",
    "18+14=",
    "34+32=",
    "17+45=",
    "60-57=",
    "52+45=",
    "56-37=",
    "99-82=",
    "59-4=",
    "19+14=",
    "25+41=",
    "68-30=",
    "50-50=",
    "2+84=",
    "72+4=",
    "31+56=",
    "58+23=",
    "28+6=",
    "78-69=",
    "12+74=",
    "22+56=",
    "16+6=",
    "43-42=",
    "85-27=",
    "97-43=",
    "61-35=",
    "62-50=",
    "14+23=",
    "63-6=",
    "83-41=",
    "60-39=",
    "48+40=",
    "1+9=",
    "40-0=",
    "51+12=",
    "86-49=",
    "60-1=",
    "10+78=",
    "77+20=",
    "91-35=",
    "65+3=",
    "44-2=",
    "40-14=",
    "66-48=",
    "56+20=",
    "72-71=",
    "58+36=",
    "35+43=",
    "13+74=",
    "46-15=",
    "42+46=",
    "31+56=",
    "88+9=",
    "30+44=",
    "24-3=",
    "34+6=",
    "63-56=",
    "36+41=",
    "31+0=",
    "72+8=",
    "75-49=",
    "86-56=",
    "44+12=",
    "15+67=",
    "8+13=",
    "61+31=",
    "85+5=",
    "4+10=",
    "38+21=",
    "85-31=",
    "85-83=",
    "2+23=",
    "52+15=",
    "99-87=",
    "81-11=",
    "15+31=",
    "80+12=",
    "54+8=",
    "53+46=",
    "87-27=",
    "22+59=",
    "35-10=",
    "47+38=",
    "78-4=",
    "10-9=",
    "39+8="
)

$idx = 0
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        if ($idx -lt $values.Length) {
            $table.Cell($r, $c).Range.Text = $values[$idx]
        }
        $idx++
    }
}

Write-Host "Done. Updated $idx cells (expected $($values.Length))."
